# Add a new "2022" column (column S) to the sheet, mirroring the layout
# of the existing "2021" column (column R): same per-cell formatting and
# one new data value (or the "-" placeholder) per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Clone the formatting of R3:R34 onto S3:S34 so every new cell picks up
#    the same number format / font / alignment / borders as its neighbour
#    in column R.
$ws.Range("R3:R34").Copy() | Out-Null
$ws.Range("S3:S34").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# 2) Fill in the 2022 values (row 4 is the year header; rows 5-34 are data,
#    with "-" used for the rows that have no figure, matching column R).
$ws.Range("S4").Value = 2022

$values = @{
    5  = 135
    6  = 99
    7  = 36
    8  = 97
    9  = 80
    10 = 17
    11 = 17
    12 = 11
    13 = 6
    14 = 5
    15 = 3
    16 = 2
    17 = "-"
    18 = "-"
    19 = "-"
    20 = 6
    21 = 1
    22 = 5
    23 = "-"
    24 = "-"
    25 = "-"
    26 = 10
    27 = 4
    28 = 6
    29 = "-"
    30 = "-"
    31 = "-"
    32 = "-"
    33 = "-"
    34 = "-"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 19).Value = $values[$row]
}

# 3) Match the saved selection recorded in the edited file.
$ws.Range("T24").Select() | Out-Null
